# Additional screenshot added in cucumber report
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SegmentManagement")
$ws.Activate()
$ws.Range("C2").Value = "12_75"
$ws.Range("C6").Select()
